$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[53.90277715644439, 73.3084717447039]"
$ws.Range("Q2").Value = "[1.4025528763774249, 1.729605565039427]"
$ws.Range("U2").Value = "[45.50448898290632, 58.63787158398123]"
$ws.Range("Y2").Value = 17.34266266266288
$ws.Range("Z2").Value = 18.5882682682685
$ws.Range("M3").Value = "[53.47426889848056, 73.9393693325495]"
$ws.Range("Q3").Value = "[1.4654476241970409, 1.7925003128590413]"
$ws.Range("U3").Value = "[43.44843140726394, 56.5515377893674]"
$ws.Range("M4").Value = "[53.4367281964297, 73.97691003460037]"
$ws.Range("N4").Value = [double]"4.440892098500626e-16"
$ws.Range("O4").Value = [double]"4.440892098500626e-16"
$ws.Range("Q4").Value = "[1.4654476241970409, 1.7925003128590422]"
$ws.Range("U4").Value = "[43.452848945511505, 56.547120251119836]"
$ws.Range("Y4").Value = 17.10312312312334
$ws.Range("Z4").Value = 18.34872872872896
$ws.Range("M5").Value = "[53.47482976993149, 73.93880846109857]"
$ws.Range("Q5").Value = "[1.46544762419704, 1.7925003128590422]"
$ws.Range("U5").Value = "[43.44694248323287, 56.553026713398474]"
$ws.Range("Y5").Value = 17.10312312312334
$ws.Range("Z5").Value = 18.34872872872896
$ws.Range("M6").Value = "[53.020956836176794, 74.19025410133956]"
$ws.Range("N6").Value = [double]"8.881784197001252e-16"
$ws.Range("O6").Value = [double]"8.881784197001252e-16"
$ws.Range("Q6").Value = "[1.503184472888809, 1.8805529598065034]"
$ws.Range("U6").Value = "[41.3608363211718, 54.496864620810285]"
$ws.Range("Y6").Value = 16.76776776776797
$ws.Range("Z6").Value = 18.20500500500524
$ws.Range("M7").Value = "[52.69438976252404, 74.78310257923046]"
$ws.Range("N7").Value = [double]"3.774758283725532e-15"
$ws.Range("O7").Value = [double]"3.774758283725532e-15"
$ws.Range("U7").Value = "[43.492933632621735, 56.55519910158505]"
$ws.Range("M8").Value = "[52.95447475186701, 74.52301758988749]"
$ws.Range("N8").Value = [double]"1.554312234475219e-15"
$ws.Range("O8").Value = [double]"1.554312234475219e-15"
$ws.Range("Q8").Value = "[1.566079220708426, 1.9434477076261185]"
$ws.Range("U8").Value = "[43.48673274675755, 56.561399987449235]"
$ws.Range("Y8").Value = 16.52822822822844
$ws.Range("Z8").Value = 17.96546546546569
$ws.Range("M9").Value = "[52.993828544227725, 74.48366379752677]"
$ws.Range("N9").Value = [double]"1.554312234475219e-15"
$ws.Range("O9").Value = [double]"1.554312234475219e-15"
$ws.Range("U9").Value = "[43.485353099167405, 56.562779635039384]"
$ws.Range("M10").Value = "[53.07467664752599, 74.18318358967745]"
$ws.Range("Q10").Value = "[1.528342372016656, 1.8553950606786564]"
$ws.Range("U10").Value = "[41.35096110356302, 54.523915081610795]"
$ws.Range("Y10").Value = 16.82834834834856
$ws.Range("Z10").Value = 18.07135135135157
$ws.Range("M11").Value = "[52.93367273589621, 74.58600686509624]"
$ws.Range("N11").Value = [double]"1.77635683940025e-15"
$ws.Range("O11").Value = [double]"1.77635683940025e-15"
$ws.Range("U11").Value = "[43.48777301167136, 56.5948749761939]"
$ws.Range("M12").Value = "[52.77192837055199, 74.74775123044046]"
$ws.Range("N12").Value = [double]"3.108624468950438e-15"
$ws.Range("O12").Value = [double]"3.108624468950438e-15"
$ws.Range("Q12").Value = "[1.57865817027235, 1.9308687580621946]"
$ws.Range("U12").Value = "[43.49079839654503, 56.591849591320226]"
$ws.Range("Y12").Value = 16.54150150150171
$ws.Range("Z12").Value = 17.88012012012033
$ws.Range("M13").Value = "[52.07891085716646, 75.50690291376527]"
$ws.Range("N13").Value = [double]"2.642330798607873e-14"
$ws.Range("O13").Value = [double]"2.642330798607873e-14"
$ws.Range("U13").Value = "[43.55572022348478, 56.60880518691177]"
$ws.Range("M14").Value = "[51.97075751521997, 75.61505625571176]"
$ws.Range("N14").Value = [double]"3.597122599785507e-14"
$ws.Range("O14").Value = [double]"3.597122599785507e-14"
$ws.Range("Q14").Value = "[1.679289766783735, 2.081816152829272]"
$ws.Range("U14").Value = "[43.54917092804926, 56.61535448234729]"
$ws.Range("Y14").Value = 15.96780780780801
$ws.Range("Z14").Value = 17.49765765765787
